$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input_ascp")

# Set in the order that matches the shared-string table insertion order
# produced by the canonical edit: DAPI, G, T, A, C
$ws.Range("G6").Value = "DAPI"
$ws.Range("G3").Value = "G"
$ws.Range("G2").Value = "T"
$ws.Range("G4").Value = "A"
$ws.Range("G5").Value = "C"

$ws.Range("A12").Select()
